$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 7860519
$ws.Range("I80").Value = 575.5625
$ws.Range("K80").Value = 1726.6875
$ws.Range("M80").Value = -728.6875
$ws.Range("H83").Value = 7860519
$ws.Range("I83").Value = 575.5625
$ws.Range("K83").Value = 5180.0625
$ws.Range("M83").Value = -188.0625
$ws.Range("H98").Value = 679.8929000000001
$ws.Range("I98").Value = 713.8333
$ws.Range("K98").Value = 713.8333
$ws.Range("M98").Value = 784.1667
$ws.Range("H122").Value = 679.8929000000001
$ws.Range("I122").Value = 713.8333
$ws.Range("K122").Value = 2141.4999
$ws.Range("M122").Value = 308.5001000000002
$ws.Range("H129").Value = 155275.89
$ws.Range("J129").Value = 157694.27
$ws.Range("L129").Value = 473082.8099999999
$ws.Range("N129").Value = -483082.8099999999
$ws.Range("H132").Value = 2455.5
$ws.Range("I132").Value = 2611.4736
$ws.Range("K132").Value = 7834.4208
$ws.Range("M132").Value = -5304.4208
$ws.Range("H138").Value = 1904.4691
$ws.Range("J138").Value = 2222.7097
$ws.Range("L138").Value = 6668.1291
$ws.Range("N138").Value = -16948.1291
$ws.Range("H140").Value = 50000
$ws.Range("J140").Value = 50000
$ws.Range("L140").Value = 50000
$ws.Range("N140").Value = -60360
$ws.Range("H141").Value = 1309.7391
$ws.Range("I141").Value = 1187.4546
$ws.Range("K141").Value = 3562.3638
$ws.Range("M141").Value = 1617.6362

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20916.2
$ws.Range("I32").Value = 24015.863
$ws.Range("K32").Value = 24015.863
$ws.Range("M32").Value = -23728.863
$ws.Range("H74").Value = 52632916
$ws.Range("I74").Value = 83333990
$ws.Range("J74").Value = 2493.5715
$ws.Range("K74").Value = 83333990
$ws.Range("L74").Value = 2493.5715
$ws.Range("M74").Value = -83333116
$ws.Range("N74").Value = -4241.5715
$ws.Range("H77").Value = 52632916
$ws.Range("I77").Value = 83333990
$ws.Range("J77").Value = 2493.5715
$ws.Range("K77").Value = 416669950
$ws.Range("L77").Value = 12467.8575
$ws.Range("M77").Value = -416665582
$ws.Range("N77").Value = -21203.8575
$ws.Range("H94").Value = 36665
$ws.Range("J94").Value = 36665
$ws.Range("L94").Value = 36665
$ws.Range("N94").Value = -38467
$ws.Range("H122").Value = 2796.9285
$ws.Range("I122").Value = 1740.8
$ws.Range("K122").Value = 5222.4
$ws.Range("M122").Value = -2772.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 26332.334
$ws.Range("J92").Value = 26332.334
$ws.Range("L92").Value = 26332.334
$ws.Range("N92").Value = -31324.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9990.927
$ws.Range("I31").Value = 18984.611
$ws.Range("J31").Value = 2952.3914
$ws.Range("K31").Value = 18984.611
$ws.Range("L31").Value = 2952.3914
$ws.Range("M31").Value = -18689.611
$ws.Range("N31").Value = -3542.3914
$ws.Range("H34").Value = 9990.927
$ws.Range("I34").Value = 18984.611
$ws.Range("J34").Value = 2952.3914
$ws.Range("K34").Value = 18984.611
$ws.Range("L34").Value = 2952.3914
$ws.Range("M34").Value = -18782.611
$ws.Range("N34").Value = -3356.3914

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 5461.1816
$ws.Range("J81").Value = 5461.1816
$ws.Range("L81").Value = 16383.5448
$ws.Range("N81").Value = -18629.5448
$ws.Range("H84").Value = 5461.1816
$ws.Range("J84").Value = 5461.1816
$ws.Range("L84").Value = 49150.6344
$ws.Range("N84").Value = -60382.6344
$ws.Range("H131").Value = 773.21
$ws.Range("I131").Value = 232.28572
$ws.Range("J131").Value = 813.92474
$ws.Range("K131").Value = 696.85716
$ws.Range("L131").Value = 2441.77422
$ws.Range("M131").Value = 4343.14284
$ws.Range("N131").Value = -12521.77422
$ws.Range("H133").Value = 3181.5
$ws.Range("I133").Value = 1536.875
$ws.Range("J133").Value = 6470.75
$ws.Range("K133").Value = 4610.625
$ws.Range("L133").Value = 19412.25
$ws.Range("M133").Value = 449.375
$ws.Range("N133").Value = -29532.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2090646.9
$ws.Range("I70").Value = 9483.833000000001
$ws.Range("J70").Value = 5212391.5
$ws.Range("K70").Value = 9483.833000000001
$ws.Range("L70").Value = 5212391.5
$ws.Range("M70").Value = -9213.833000000001
$ws.Range("N70").Value = -5212931.5
$ws.Range("H73").Value = 2090646.9
$ws.Range("I73").Value = 9483.833000000001
$ws.Range("J73").Value = 5212391.5
$ws.Range("K73").Value = 9483.833000000001
$ws.Range("L73").Value = 5212391.5
$ws.Range("M73").Value = -8547.833000000001
$ws.Range("N73").Value = -5214263.5
$ws.Range("H102").Value = 35715680
$ws.Range("I102").Value = 41668000
$ws.Range("K102").Value = 41668000
$ws.Range("M102").Value = -41666378
$ws.Range("H122").Value = 60607944
$ws.Range("I122").Value = 23811056
$ws.Range("J122").Value = 125002500
$ws.Range("K122").Value = 71433168
$ws.Range("L122").Value = 375007500
$ws.Range("M122").Value = -71430718
$ws.Range("N122").Value = -375012400
$ws.Range("H126").Value = 5206.9165
$ws.Range("I126").Value = 3956.25
$ws.Range("J126").Value = 7708.25
$ws.Range("K126").Value = 11868.75
$ws.Range("L126").Value = 23124.75
$ws.Range("M126").Value = -9398.75
$ws.Range("N126").Value = -28064.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4128.615
$ws.Range("I7").Value = 2362.1667
$ws.Range("K7").Value = 2362.1667
$ws.Range("M7").Value = -2250.1667
$ws.Range("H40").Value = 6635.5
$ws.Range("I40").Value = 4750
$ws.Range("J40").Value = 7892.5
$ws.Range("K40").Value = 4750
$ws.Range("L40").Value = 7892.5
$ws.Range("M40").Value = -4614
$ws.Range("N40").Value = -8164.5
$ws.Range("H122").Value = 1403181.5
$ws.Range("I122").Value = 2180982.2
$ws.Range("K122").Value = 6542946.600000001
$ws.Range("M122").Value = -6540496.600000001
$ws.Range("H126").Value = 4128.615
$ws.Range("I126").Value = 2362.1667
$ws.Range("K126").Value = 7086.500100000001
$ws.Range("M126").Value = -4616.500100000001
$ws.Range("H132").Value = 2442.9375
$ws.Range("I132").Value = 1822.1111
$ws.Range("J132").Value = 3241.1428
$ws.Range("K132").Value = 5466.3333
$ws.Range("L132").Value = 9723.428400000001
$ws.Range("M132").Value = -2936.3333
$ws.Range("N132").Value = -14783.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2192.4285
$ws.Range("I122").Value = 2058
$ws.Range("K122").Value = 6174
$ws.Range("M122").Value = -3724
$ws.Range("H132").Value = 1018
$ws.Range("I132").Value = 599.375
$ws.Range("J132").Value = 2357.6
$ws.Range("K132").Value = 1798.125
$ws.Range("L132").Value = 7072.799999999999
$ws.Range("M132").Value = 731.875
$ws.Range("N132").Value = -12132.8

Write-Output "Applied all updates"